## Scheduled market-data refresh: update Leve profit calculations
## (currentAveragePrice* / LevePrice* / LeveProfit* columns) per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 64456.867
$ws.Cells.Item(74, 9).Value = 116481.625
$ws.Cells.Item(74, 11).Value = 116481.625
$ws.Cells.Item(74, 13).Value = -115545.625
$ws.Cells.Item(77, 8).Value = 64456.867
$ws.Cells.Item(77, 9).Value = 116481.625
$ws.Cells.Item(77, 11).Value = 582408.125
$ws.Cells.Item(77, 13).Value = -577728.125
$ws.Cells.Item(97, 8).Value = 1582.7273
$ws.Cells.Item(97, 9).Value = 533.3333
$ws.Cells.Item(97, 10).Value = 1976.25
$ws.Cells.Item(97, 11).Value = 1599.9999
$ws.Cells.Item(97, 12).Value = 5928.75
$ws.Cells.Item(97, 13).Value = -1103.9999
$ws.Cells.Item(97, 14).Value = -6920.75
$ws.Cells.Item(106, 8).Value = 3826.25
$ws.Cells.Item(106, 9).Value = 5752.5
$ws.Cells.Item(106, 10).Value = 1900
$ws.Cells.Item(106, 11).Value = 5752.5
$ws.Cells.Item(106, 12).Value = 1900
$ws.Cells.Item(106, 13).Value = -5121.5
$ws.Cells.Item(106, 14).Value = -3162
$ws.Cells.Item(132, 8).Value = 2042.3334
$ws.Cells.Item(132, 9).Value = 1714.0714
$ws.Cells.Item(132, 11).Value = 5142.2142
$ws.Cells.Item(132, 13).Value = -2612.2142
$ws.Cells.Item(134, 8).Value = 77816.664
$ws.Cells.Item(134, 10).Value = 77816.664
$ws.Cells.Item(134, 12).Value = 77816.664
$ws.Cells.Item(134, 14).Value = -87956.664
$ws.Cells.Item(138, 8).Value = 3729.6296
$ws.Cells.Item(138, 9).Value = 1425.3914
$ws.Cells.Item(138, 10).Value = 5439.2256
$ws.Cells.Item(138, 11).Value = 4276.174199999999
$ws.Cells.Item(138, 12).Value = 16317.6768
$ws.Cells.Item(138, 13).Value = 863.8258000000005
$ws.Cells.Item(138, 14).Value = -26597.6768
$ws.Cells.Item(141, 8).Value = 5169.853
$ws.Cells.Item(141, 9).Value = 2179.516
$ws.Cells.Item(141, 10).Value = 36070
$ws.Cells.Item(141, 11).Value = 6538.548000000001
$ws.Cells.Item(141, 12).Value = 108210
$ws.Cells.Item(141, 13).Value = -1358.548000000001
$ws.Cells.Item(141, 14).Value = -118570

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 26500
$ws.Cells.Item(3, 9).Value = 26500
$ws.Cells.Item(3, 11).Value = 26500
$ws.Cells.Item(3, 13).Value = -26385
$ws.Cells.Item(5, 8).Value = 260
$ws.Cells.Item(5, 9).Value = 270
$ws.Cells.Item(5, 11).Value = 270
$ws.Cells.Item(5, 13).Value = -158
$ws.Cells.Item(22, 8).Value = 987.5
$ws.Cells.Item(22, 9).Value = 987.5
$ws.Cells.Item(22, 11).Value = 987.5
$ws.Cells.Item(22, 13).Value = -688.5
$ws.Cells.Item(32, 8).Value = 14470.807
$ws.Cells.Item(32, 9).Value = 14534.078
$ws.Cells.Item(32, 10).Value = 13933
$ws.Cells.Item(32, 11).Value = 14534.078
$ws.Cells.Item(32, 12).Value = 13933
$ws.Cells.Item(32, 13).Value = -14247.078
$ws.Cells.Item(32, 14).Value = -14507
$ws.Cells.Item(61, 8).Value = 3481.6667
$ws.Cells.Item(61, 9).Value = 3951.6
$ws.Cells.Item(61, 10).Value = 3205.2354
$ws.Cells.Item(61, 11).Value = 3951.6
$ws.Cells.Item(61, 12).Value = 3205.2354
$ws.Cells.Item(61, 13).Value = -3739.6
$ws.Cells.Item(61, 14).Value = -3629.2354
$ws.Cells.Item(63, 8).Value = 4623.75
$ws.Cells.Item(63, 9).Value = 3798
$ws.Cells.Item(63, 11).Value = 3798
$ws.Cells.Item(63, 13).Value = -3112
$ws.Cells.Item(66, 8).Value = 4623.75
$ws.Cells.Item(66, 9).Value = 3798
$ws.Cells.Item(66, 11).Value = 18990
$ws.Cells.Item(66, 13).Value = -15558
$ws.Cells.Item(74, 8).Value = 2661.3438
$ws.Cells.Item(74, 9).Value = 3223.6875
$ws.Cells.Item(74, 10).Value = 2099
$ws.Cells.Item(74, 11).Value = 3223.6875
$ws.Cells.Item(74, 12).Value = 2099
$ws.Cells.Item(74, 13).Value = -2349.6875
$ws.Cells.Item(74, 14).Value = -3847
$ws.Cells.Item(77, 8).Value = 2661.3438
$ws.Cells.Item(77, 9).Value = 3223.6875
$ws.Cells.Item(77, 10).Value = 2099
$ws.Cells.Item(77, 11).Value = 16118.4375
$ws.Cells.Item(77, 12).Value = 10495
$ws.Cells.Item(77, 13).Value = -11750.4375
$ws.Cells.Item(77, 14).Value = -19231
$ws.Cells.Item(88, 8).Value = 2650.8333
$ws.Cells.Item(88, 9).Value = 1937.3334
$ws.Cells.Item(88, 10).Value = 2888.6667
$ws.Cells.Item(88, 11).Value = 1937.3334
$ws.Cells.Item(88, 12).Value = 2888.6667
$ws.Cells.Item(88, 13).Value = -1531.3334
$ws.Cells.Item(88, 14).Value = -3700.6667
$ws.Cells.Item(91, 8).Value = 2650.8333
$ws.Cells.Item(91, 9).Value = 1937.3334
$ws.Cells.Item(91, 10).Value = 2888.6667
$ws.Cells.Item(91, 11).Value = 1937.3334
$ws.Cells.Item(91, 12).Value = 2888.6667
$ws.Cells.Item(91, 13).Value = -533.3334
$ws.Cells.Item(91, 14).Value = -5696.6667
$ws.Cells.Item(132, 8).Value = 5323.9023
$ws.Cells.Item(132, 9).Value = 8536.588
$ws.Cells.Item(132, 10).Value = 3048.25
$ws.Cells.Item(132, 11).Value = 25609.764
$ws.Cells.Item(132, 12).Value = 9144.75
$ws.Cells.Item(132, 13).Value = -23079.764
$ws.Cells.Item(132, 14).Value = -14204.75
$ws.Cells.Item(136, 8).Value = 3481.6667
$ws.Cells.Item(136, 9).Value = 3951.6
$ws.Cells.Item(136, 10).Value = 3205.2354
$ws.Cells.Item(136, 11).Value = 11854.8
$ws.Cells.Item(136, 12).Value = 9615.706200000001
$ws.Cells.Item(136, 13).Value = -9304.799999999999
$ws.Cells.Item(136, 14).Value = -14715.7062

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 260
$ws.Cells.Item(4, 9).Value = 270
$ws.Cells.Item(4, 11).Value = 270
$ws.Cells.Item(4, 13).Value = -155
$ws.Cells.Item(35, 8).Value = 39741.332
$ws.Cells.Item(35, 10).Value = 39741.332
$ws.Cells.Item(35, 12).Value = 39741.332
$ws.Cells.Item(35, 14).Value = -40361.332
$ws.Cells.Item(82, 8).Value = 41740
$ws.Cells.Item(82, 9).Value = 41740
$ws.Cells.Item(82, 11).Value = 41740
$ws.Cells.Item(82, 13).Value = -41357
$ws.Cells.Item(85, 8).Value = 41740
$ws.Cells.Item(85, 9).Value = 41740
$ws.Cells.Item(85, 11).Value = 41740
$ws.Cells.Item(85, 13).Value = -40414

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3893.4583
$ws.Cells.Item(31, 9).Value = 1615.8948
$ws.Cells.Item(31, 10).Value = 5385.6553
$ws.Cells.Item(31, 11).Value = 1615.8948
$ws.Cells.Item(31, 12).Value = 5385.6553
$ws.Cells.Item(31, 13).Value = -1320.8948
$ws.Cells.Item(31, 14).Value = -5975.6553
$ws.Cells.Item(34, 8).Value = 3893.4583
$ws.Cells.Item(34, 9).Value = 1615.8948
$ws.Cells.Item(34, 10).Value = 5385.6553
$ws.Cells.Item(34, 11).Value = 1615.8948
$ws.Cells.Item(34, 12).Value = 5385.6553
$ws.Cells.Item(34, 13).Value = -1413.8948
$ws.Cells.Item(34, 14).Value = -5789.6553
$ws.Cells.Item(58, 8).Value = 1858.525
$ws.Cells.Item(58, 9).Value = 1573.3448
$ws.Cells.Item(58, 10).Value = 2610.3635
$ws.Cells.Item(58, 11).Value = 1573.3448
$ws.Cells.Item(58, 12).Value = 2610.3635
$ws.Cells.Item(58, 13).Value = -1370.3448
$ws.Cells.Item(58, 14).Value = -3016.3635
$ws.Cells.Item(122, 8).Value = 1578.5294
$ws.Cells.Item(122, 9).Value = 1571.9231
$ws.Cells.Item(122, 11).Value = 4715.7693
$ws.Cells.Item(122, 13).Value = -2265.7693
$ws.Cells.Item(132, 8).Value = 2141.0386
$ws.Cells.Item(132, 9).Value = 1507.75
$ws.Cells.Item(132, 10).Value = 4252
$ws.Cells.Item(132, 11).Value = 4523.25
$ws.Cells.Item(132, 12).Value = 12756
$ws.Cells.Item(132, 13).Value = -1993.25
$ws.Cells.Item(132, 14).Value = -17816
$ws.Cells.Item(136, 8).Value = 1858.525
$ws.Cells.Item(136, 9).Value = 1573.3448
$ws.Cells.Item(136, 10).Value = 2610.3635
$ws.Cells.Item(136, 11).Value = 4720.0344
$ws.Cells.Item(136, 12).Value = 7831.0905
$ws.Cells.Item(136, 13).Value = -2170.0344
$ws.Cells.Item(136, 14).Value = -12931.0905
$ws.Cells.Item(141, 8).Value = 25000
$ws.Cells.Item(141, 10).Value = 25000
$ws.Cells.Item(141, 12).Value = 25000
$ws.Cells.Item(141, 14).Value = -35360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 8191.3335
$ws.Cells.Item(17, 10).Value = 8117.8184
$ws.Cells.Item(17, 12).Value = 24353.4552
$ws.Cells.Item(17, 14).Value = -24691.4552
$ws.Cells.Item(34, 8).Value = 1271.9333
$ws.Cells.Item(34, 10).Value = 2896.5
$ws.Cells.Item(34, 12).Value = 8689.5
$ws.Cells.Item(34, 14).Value = -8857.5
$ws.Cells.Item(39, 8).Value = 4524.467
$ws.Cells.Item(39, 10).Value = 4524.467
$ws.Cells.Item(39, 12).Value = 13573.401
$ws.Cells.Item(39, 14).Value = -14161.401
$ws.Cells.Item(55, 8).Value = 5953.909
$ws.Cells.Item(55, 10).Value = 6449.3
$ws.Cells.Item(55, 12).Value = 19347.9
$ws.Cells.Item(55, 14).Value = -19701.9
$ws.Cells.Item(133, 8).Value = 3924.577
$ws.Cells.Item(133, 9).Value = 1939.25
$ws.Cells.Item(133, 10).Value = 5626.2856
$ws.Cells.Item(133, 11).Value = 5817.75
$ws.Cells.Item(133, 12).Value = 16878.8568
$ws.Cells.Item(133, 13).Value = -757.75
$ws.Cells.Item(133, 14).Value = -26998.8568
$ws.Cells.Item(139, 8).Value = 1615.7838
$ws.Cells.Item(139, 9).Value = 1104.3704
$ws.Cells.Item(139, 11).Value = 3313.1112
$ws.Cells.Item(139, 13).Value = 1826.8888

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3256
$ws.Cells.Item(132, 9).Value = 2390.2856
$ws.Cells.Item(132, 11).Value = 7170.8568
$ws.Cells.Item(132, 13).Value = -4640.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(133, 8).Value = 39250
$ws.Cells.Item(133, 10).Value = 39250
$ws.Cells.Item(133, 12).Value = 39250
$ws.Cells.Item(133, 14).Value = -44310
$ws.Cells.Item(136, 8).Value = 3100.2222
$ws.Cells.Item(136, 9).Value = 2317
$ws.Cells.Item(136, 10).Value = 4666.6665
$ws.Cells.Item(136, 11).Value = 6951
$ws.Cells.Item(136, 12).Value = 13999.9995
$ws.Cells.Item(136, 13).Value = -4401
$ws.Cells.Item(136, 14).Value = -19099.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 13157980
$ws.Cells.Item(39, 10).Value = 40000
$ws.Cells.Item(39, 12).Value = 40000
$ws.Cells.Item(39, 14).Value = -40826
$ws.Cells.Item(132, 8).Value = 4220.174
$ws.Cells.Item(132, 9).Value = 4237.8
$ws.Cells.Item(132, 10).Value = 4187.125
$ws.Cells.Item(132, 11).Value = 12713.4
$ws.Cells.Item(132, 12).Value = 12561.375
$ws.Cells.Item(132, 13).Value = -10183.4
$ws.Cells.Item(132, 14).Value = -17621.375
